$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50 - this shifts rows 50:90 down to 51:91
# and keeps the existing dimension/data intact.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new data record.
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 45159
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 100112031
$ws.Range("G50").Value = "Poroto verde"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 2200
$ws.Range("K50").Value = 1300
$ws.Range("L50").Value = 1400
$ws.Range("M50").Value = 1345
$ws.Range("N50").Value = "`$/kilo"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value = 1345
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"
